# "Actualització de les hores" — add a new log row (date + hours) to the
# "Manel" sheet, mirroring the existing row-1 pattern (date in col B styled
# as a date, count in col C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Manel")

# New row 2: B2 = date (11/11/2025, serial 45972), C2 = 1.
# Copy B1's formatting (date number format) down to B2 first, then set the
# values so the new row matches the existing date-styled column.
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B2").Value = 45972

$ws.Range("C2").Value = 1

# Move the active selection to C2, matching the post-edit state.
$ws.Range("C2").Select()
